# Scheduled-runner refresh of leve-profit market data (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the crafting-class sheets. Updates the
# H:N columns for the affected leve rows on each sheet; where a profit
# column had no prior value it is newly populated, and where a column no
# longer applies its cell is cleared back out entirely (not zeroed).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 695
$ws.Range("I61").Value = 695
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2085
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1913
$ws.Range("N61").Value = ""
$ws.Range("H70").Value = 29735686
$ws.Range("J70").Value = 38467880
$ws.Range("L70").Value = 115403640
$ws.Range("N70").Value = -115404180
$ws.Range("H73").Value = 29735686
$ws.Range("J73").Value = 38467880
$ws.Range("L73").Value = 115403640
$ws.Range("N73").Value = -115405512
$ws.Range("H86").Value = 3270.111
$ws.Range("J86").Value = 3186.4
$ws.Range("L86").Value = 3186.4
$ws.Range("N86").Value = -5432.4
$ws.Range("H88").Value = 3417.9524
$ws.Range("I88").Value = 699.5
$ws.Range("J88").Value = 3704.1052
$ws.Range("K88").Value = 699.5
$ws.Range("L88").Value = 3704.1052
$ws.Range("M88").Value = -293.5
$ws.Range("N88").Value = -4516.1052
$ws.Range("H89").Value = 3270.111
$ws.Range("J89").Value = 3186.4
$ws.Range("L89").Value = 15932
$ws.Range("N89").Value = -27164
$ws.Range("H91").Value = 3417.9524
$ws.Range("I91").Value = 699.5
$ws.Range("J91").Value = 3704.1052
$ws.Range("K91").Value = 699.5
$ws.Range("L91").Value = 3704.1052
$ws.Range("M91").Value = 704.5
$ws.Range("N91").Value = -6512.1052
$ws.Range("H132").Value = 1398.9592
$ws.Range("J132").Value = 2146.6667
$ws.Range("L132").Value = 6440.000100000001
$ws.Range("N132").Value = -11500.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = ""
$ws.Range("H50").Value = 8852.444
$ws.Range("I50").Value = 673.3333
$ws.Range("J50").Value = 12942
$ws.Range("K50").Value = 673.3333
$ws.Range("L50").Value = 12942
$ws.Range("M50").Value = 40.66669999999999
$ws.Range("N50").Value = -14370

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = ""
$ws.Range("H86").Value = 2089
$ws.Range("I86").Value = 2089
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2089
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -966
$ws.Range("N86").Value = ""
$ws.Range("H89").Value = 2089
$ws.Range("I89").Value = 2089
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10445
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4829
$ws.Range("N89").Value = ""
$ws.Range("H99").Value = 3336.5715
$ws.Range("I99").Value = 3154.4666
$ws.Range("K99").Value = 3154.4666
$ws.Range("M99").Value = -1656.4666
$ws.Range("H105").Value = 2218.3547
$ws.Range("I105").Value = 2243.5
$ws.Range("J105").Value = 2132.1428
$ws.Range("K105").Value = 2243.5
$ws.Range("L105").Value = 2132.1428
$ws.Range("M105").Value = -496.5
$ws.Range("N105").Value = -5626.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4833.3335
$ws.Range("I16").Value = 4750
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 4750
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -4463
$ws.Range("N16").Value = -5574
$ws.Range("H31").Value = 5705.7666
$ws.Range("I31").Value = 2888.9473
$ws.Range("J31").Value = 10571.182
$ws.Range("K31").Value = 2888.9473
$ws.Range("L31").Value = 10571.182
$ws.Range("M31").Value = -2593.9473
$ws.Range("N31").Value = -11161.182
$ws.Range("H34").Value = 5705.7666
$ws.Range("I34").Value = 2888.9473
$ws.Range("J34").Value = 10571.182
$ws.Range("K34").Value = 2888.9473
$ws.Range("L34").Value = 10571.182
$ws.Range("M34").Value = -2686.9473
$ws.Range("N34").Value = -10975.182
$ws.Range("H50").Value = 46664
$ws.Range("J50").Value = 46664
$ws.Range("L50").Value = 46664
$ws.Range("N50").Value = -47914
$ws.Range("H51").Value = 33998.668
$ws.Range("J51").Value = 33998.668
$ws.Range("L51").Value = 33998.668
$ws.Range("N51").Value = -35470.668
$ws.Range("H56").Value = 32499.5
$ws.Range("J56").Value = 49999
$ws.Range("L56").Value = 49999
$ws.Range("N56").Value = -51689
$ws.Range("H58").Value = 2300.1333
$ws.Range("I58").Value = 2250.3333
$ws.Range("J58").Value = 2499.3333
$ws.Range("K58").Value = 2250.3333
$ws.Range("L58").Value = 2499.3333
$ws.Range("M58").Value = -2047.3333
$ws.Range("N58").Value = -2905.3333
$ws.Range("H59").Value = 41243
$ws.Range("J59").Value = 45114.5
$ws.Range("L59").Value = 45114.5
$ws.Range("N59").Value = -47404.5
$ws.Range("H60").Value = 39999
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").Value = ""
$ws.Range("H61").Value = 33998.668
$ws.Range("J61").Value = 33998.668
$ws.Range("L61").Value = 33998.668
$ws.Range("N61").Value = -34694.668
$ws.Range("H68").Value = 75000
$ws.Range("J68").Value = 75000
$ws.Range("L68").Value = 75000
$ws.Range("N68").Value = -76498
$ws.Range("H71").Value = 75000
$ws.Range("J71").Value = 75000
$ws.Range("L71").Value = 225000
$ws.Range("N71").Value = -232488
$ws.Range("J74").Value = 45000
$ws.Range("L74").Value = 45000
$ws.Range("N74").Value = -46748
$ws.Range("J77").Value = 45000
$ws.Range("L77").Value = 135000
$ws.Range("N77").Value = -143736
$ws.Range("H105").Value = 1171.5
$ws.Range("J105").Value = 1011.6667
$ws.Range("L105").Value = 1011.6667
$ws.Range("N105").Value = -4505.6667
$ws.Range("H113").Value = 4833.3335
$ws.Range("I113").Value = 4750
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 4750
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -2580
$ws.Range("N113").Value = -9340
$ws.Range("H122").Value = 3619.9607
$ws.Range("I122").Value = 765.1
$ws.Range("K122").Value = 2295.3
$ws.Range("M122").Value = 154.6999999999998
$ws.Range("H134").Value = 3098.1365
$ws.Range("I134").Value = 3352.5881
$ws.Range("K134").Value = 10057.7643
$ws.Range("M134").Value = -7522.764299999999
$ws.Range("H136").Value = 2300.1333
$ws.Range("I136").Value = 2250.3333
$ws.Range("J136").Value = 2499.3333
$ws.Range("K136").Value = 6750.999899999999
$ws.Range("L136").Value = 7497.999899999999
$ws.Range("M136").Value = -4200.999899999999
$ws.Range("N136").Value = -12597.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 197.33333
$ws.Range("I2").Value = 298.2
$ws.Range("K2").Value = 1789.2
$ws.Range("M2").Value = -1676.2
$ws.Range("H74").Value = 17432.8
$ws.Range("I74").Value = 17505.6
$ws.Range("J74").Value = 17360
$ws.Range("K74").Value = 52516.8
$ws.Range("L74").Value = 52080
$ws.Range("M74").Value = -51455.8
$ws.Range("N74").Value = -54202
$ws.Range("H77").Value = 17432.8
$ws.Range("I77").Value = 17505.6
$ws.Range("J77").Value = 17360
$ws.Range("K77").Value = 157550.4
$ws.Range("L77").Value = 156240
$ws.Range("M77").Value = -152246.4
$ws.Range("N77").Value = -166848

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1180.7
$ws.Range("I97").Value = 1066.3334
$ws.Range("K97").Value = 1066.3334
$ws.Range("M97").Value = -570.3334
$ws.Range("H132").Value = 3544.2
$ws.Range("I132").Value = 3346.1765
$ws.Range("K132").Value = 10038.5295
$ws.Range("M132").Value = -7508.529500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 25666.666
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").Value = ""
$ws.Range("H40").Value = 5796.5586
$ws.Range("I40").Value = 4421.722
$ws.Range("K40").Value = 4421.722
$ws.Range("M40").Value = -4285.722
$ws.Range("H46").Value = 3014.2
$ws.Range("I46").Value = 1995.5714
$ws.Range("J46").Value = 3905.5
$ws.Range("K46").Value = 1995.5714
$ws.Range("L46").Value = 3905.5
$ws.Range("M46").Value = -1807.5714
$ws.Range("N46").Value = -4281.5
$ws.Range("H55").Value = 2387.1667
$ws.Range("I55").Value = 1775
$ws.Range("J55").Value = 2999.3333
$ws.Range("K55").Value = 1775
$ws.Range("L55").Value = 2999.3333
$ws.Range("M55").Value = -1602
$ws.Range("N55").Value = -3345.3333
$ws.Range("H59").Value = 155750
$ws.Range("J59").Value = 155750
$ws.Range("L59").Value = 155750
$ws.Range("N59").Value = -157058
$ws.Range("H61").Value = 4514.524
$ws.Range("I61").Value = 2369.6924
$ws.Range("K61").Value = 2369.6924
$ws.Range("M61").Value = -2167.6924
$ws.Range("H82").Value = 4188.9
$ws.Range("I82").Value = 2222.5
$ws.Range("J82").Value = 5499.8335
$ws.Range("K82").Value = 2222.5
$ws.Range("L82").Value = 5499.8335
$ws.Range("M82").Value = -1861.5
$ws.Range("N82").Value = -6221.8335
$ws.Range("H85").Value = 4188.9
$ws.Range("I85").Value = 2222.5
$ws.Range("J85").Value = 5499.8335
$ws.Range("K85").Value = 2222.5
$ws.Range("L85").Value = 5499.8335
$ws.Range("M85").Value = -974.5
$ws.Range("N85").Value = -7995.8335
$ws.Range("H113").Value = 4514.524
$ws.Range("I113").Value = 2369.6924
$ws.Range("K113").Value = 2369.6924
$ws.Range("M113").Value = -199.6923999999999
$ws.Range("H140").Value = 99429
$ws.Range("J140").Value = 99429
$ws.Range("L140").Value = 99429
$ws.Range("N140").Value = -109789

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 469.16666
$ws.Range("I4").Value = 397.14285
$ws.Range("J4").Value = 570
$ws.Range("K4").Value = 397.14285
$ws.Range("L4").Value = 570
$ws.Range("M4").Value = -284.14285
$ws.Range("N4").Value = -796
$ws.Range("H49").Value = 31314.057
$ws.Range("I49").Value = 30551.45
$ws.Range("K49").Value = 30551.45
$ws.Range("M49").Value = -30321.45
